$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: path EFF_MEAN ~ X
$ws.Range("B2").Value = -0.3930257200655047
$ws.Range("C2").Value = 0.1465453386822315
$ws.Range("D2").Value = 0.008757824312620524
$ws.Range("E2").Value = -0.6843004488585831
$ws.Range("F2").Value = -0.1017509912724263

# Row 3: Y ~ EFF_MEAN
$ws.Range("B3").Value = 0.5058626820182867
$ws.Range("C3").Value = 0.06755561790893917
$ws.Range("D3").Value = 0.0000000000536522182074112
$ws.Range("E3").Value = 0.3715885763946122
$ws.Range("F3").Value = 0.640136787641961

# Row 4: Total
$ws.Range("B4").Value = -0.01360658896060151
$ws.Range("C4").Value = 0.1232051502995288
$ws.Range("D4").Value = 0.9123160068126703
$ws.Range("E4").Value = -0.2584901660303711
$ws.Range("F4").Value = 0.2312769881091681

# Row 5: Direct
$ws.Range("B5").Value = 0.2005228879382159
$ws.Range("C5").Value = 0.09820222365217117
$ws.Range("D5").Value = 0.04422159303169759
$ws.Range("E5").Value = 0.005303328411330299
$ws.Range("F5").Value = 0.3957424474651016
$ws.Range("G5").Value = "Yes"

# Row 6: Indirect
$ws.Range("B6").Value = -0.2141294768988164
$ws.Range("C6").Value = 0.07379528984726182
$ws.Range("D6").Value = 0.0012
$ws.Range("E6").Value = -0.3763093020778983
$ws.Range("F6").Value = -0.0839923724484373
